$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so values that look numeric
# (e.g. "1.01") are stored as text, matching the source data which is all inline strings.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '90.344.32'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '3.153.81'
$ws.Range('E3').Value = '  -3.87%  '
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +0.71%  '
$ws.Range('D5').Value = '212.98'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('D6').Value = '624.05'
$ws.Range('E6').Value = '  -0.55%  '
$ws.Range('D7').Value = '0.398'
$ws.Range('E7').Value = '  -3.64%  '
$ws.Range('D8').Value = '0.723'
$ws.Range('E8').Value = '  +1.73%  '
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.49%  '
$ws.Range('D10').Value = '3.151.76'
$ws.Range('E10').Value = '  -3.17%  '
$ws.Range('D11').Value = '0.558'
$ws.Range('E11').Value = '  -6.18%  '
$ws.Range('D12').Value = '0.182'
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('D13').Value = '0.0000254'
$ws.Range('E13').Value = '  -7.70%  '
$ws.Range('D14').Value = '90.160.83'
$ws.Range('E14').Value = '  -0.05%  '
$ws.Range('D15').Value = '5.27'
$ws.Range('E15').Value = '  -3.00%  '
$ws.Range('D16').Value = '3.748.02'
$ws.Range('E16').Value = '  -3.18%  '
$ws.Range('D17').Value = '32.05'
$ws.Range('E17').Value = '  -6.51%  '
$ws.Range('D18').Value = '3.172.52'
$ws.Range('E18').Value = '  -2.19%  '
$ws.Range('D19').Value = '3.30'
$ws.Range('E19').Value = '  +1.56%  '
$ws.Range('D20').Value = '0.0000214'
$ws.Range('E20').Value = '  +12.05%  '
$ws.Range('D21').Value = '13.22'
$ws.Range('E21').Value = '  -6.97%  '
$ws.Range('D22').Value = '425.83'
$ws.Range('E22').Value = '  -1.65%  '
$ws.Range('D23').Value = '8.39'
$ws.Range('E23').Value = '  -6.66%  '
$ws.Range('D24').Value = '4.90'
$ws.Range('E24').Value = '  -8.29%  '
$ws.Range('D25').Value = '5.24'
$ws.Range('E25').Value = '  -4.02%  '
$ws.Range('D26').Value = '11.48'
$ws.Range('E26').Value = '  -4.79%  '
$ws.Range('D27').Value = '79.93'
$ws.Range('E27').Value = '  +5.12%  '
$ws.Range('D28').Value = '3.351.82'
$ws.Range('E28').Value = '  -1.71%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('B31').Value = 'Cronos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D31').Value = '0.157'
$ws.Range('E31').Value = '  -11.00%  '
$ws.Range('D32').Value = '3.97'
$ws.Range('E32').Value = '  +11.41%  '
$ws.Range('D33').Value = '8.24'
$ws.Range('E33').Value = '  -4.69%  '
$ws.Range('D34').Value = '507.91'
$ws.Range('E34').Value = '  -10.39%  '
$ws.Range('D35').Value = '6.82'
$ws.Range('E35').Value = '  -5.54%  '
$ws.Range('D36').Value = '1.86'
$ws.Range('E36').Value = '  -2.91%  '
$ws.Range('D37').Value = '1.26'
$ws.Range('E37').Value = '  -8.20%  '
$ws.Range('D38').Value = '22.09'
$ws.Range('E38').Value = '  -2.74%  '
$ws.Range('D39').Value = '22.33'
$ws.Range('E39').Value = '  -0.20%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value = '1.01'
$ws.Range('E40').Value = '  +1.12%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.126'
$ws.Range('E41').Value = '  -5.79%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = '1.88'
$ws.Range('E43').Value = '  -5.39%  '
$ws.Range('D44').Value = '0.366'
$ws.Range('E44').Value = '  -7.24%  '
$ws.Range('D45').Value = '147.14'
$ws.Range('E45').Value = '  -1.23%  '
$ws.Range('D46').Value = '43.91'
$ws.Range('E46').Value = '  -0.13%  '
$ws.Range('D47').Value = '167.18'
$ws.Range('E47').Value = '  -8.97%  '
$ws.Range('D48').Value = '0.125'
$ws.Range('E48').Value = '  -3.51%  '
$ws.Range('D49').Value = '0.730'
$ws.Range('E49').Value = '  +2.05%  '
$ws.Range('D50').Value = '24.40'
$ws.Range('E50').Value = '  -3.39%  '
$ws.Range('D51').Value = '1.19'
$ws.Range('E51').Value = '  -7.96%  '

# Reset to the default style so no stray number-format style is left on the cells
# (matches the original workbook where these cells carry no explicit style index).
$ws.Range("D2:E51").Style = "Normal"
